$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: row index numbers (0-based), rows 2..24 -> values 0..22
# Rows 2..14 already exist with the correct style; extend the same style
# (bold, bordered, centered) down to the newly added rows 15..24.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A15:A24").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -le 22; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# Column B: Buying Opportunity (rows 2..14), rows 15..24 blank
$colB = @(
    "NSE:AHLEAST",
    "NSE:APOLLOHOSP",
    "NSE:DALBHARAT",
    "NSE:GODREJPROP",
    "NSE:INDTERRAIN",
    "NSE:JKCEMENT",
    "NSE:JUBLINGREA",
    "NSE:KRSNAA",
    "NSE:LTTS",
    "NSE:MANORAMA",
    "NSE:MANYAVAR",
    "NSE:OIL",
    "NSE:PHOENIXLTD"
)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 2).Value = ""
}

# Column C: support Zone -> all cleared (blank) for rows 2..24
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Value = ""
}

# Column D: long buildup (rows 2..24), 23 values
$colD = @(
    "NSE:APOLLOTYRE",
    "NSE:BAJAJFINSV",
    "NSE:BEL",
    "NSE:BHARATFORG",
    "NSE:BHARTIARTL",
    "NSE:BHEL",
    "NSE:BRITANNIA",
    "NSE:BSOFT",
    "NSE:CIPLA",
    "NSE:DABUR",
    "NSE:DEEPAKNTR",
    "NSE:GODREJPROP",
    "NSE:HEROMOTOCO",
    "NSE:INDIAMART",
    "NSE:IPCALAB",
    "NSE:IRCTC",
    "NSE:LT",
    "NSE:M&M",
    "NSE:MARICO",
    "NSE:MUTHOOTFIN",
    "NSE:NESTLEIND",
    "NSE:PEL",
    "NSE:RECLTD"
)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

# Column E: Short buildup -> all cleared (blank) for rows 2..24
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 5).Value = ""
}

# Column F: FII ENTERING (rows 2..14), rows 15..24 blank
$colF = @(
    "NSE:APOLLOHOSP",
    "NSE:COFORGE",
    "NSE:DLF",
    "NSE:GODREJPROP",
    "NSE:GRASIM",
    "NSE:HAVELLS",
    "NSE:ICICIBANK",
    "NSE:INDIGO",
    "NSE:KOTAKBANK",
    "NSE:LAURUSLABS",
    "NSE:LTTS",
    "NSE:LUPIN",
    "NSE:PERSISTENT"
)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}
for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 6).Value = ""
}
